# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Mon Sep 18 14:11:33 UTC 2023 with GitHub Actions"
#
# Updates Price (D) / Volume(1h) (E) text values for rows 2-51, and swaps the
# Algorand / EnergySwap rows (49 <-> 50, including B/C/D/E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of these columns hold text (not numeric) values in the source data,
# e.g. "27.337.66" and "  +2.05%  " -- force Text format so Excel doesn't
# coerce them (and drop formatting like trailing zeros) when parsed as numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.305.35'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.660.17'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.21%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.57'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.504'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.30%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.81'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0848'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.64%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.895.60'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.655.85'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.35%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.21%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.87'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.305.36'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0735'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '222.33'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.94%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +8.72%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.50'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +5.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.24'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.42'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.39%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.43'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.62%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.86%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0514'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.26%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.39'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.01'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.56'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.259.76'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.74%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.34%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.537'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.828'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.41%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.34%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.807.05'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.13'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.70'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '92.66'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.97%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.61%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.72'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0982'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.83%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.48%  '
